# ---------------------------------------------------------------------------
# Applies the "EarlyRePayment-Makerepayment1" edit:
#   - Recalculated repayment figures on Summary / Repayment Schedule /
#     Transactions sheets (a new, slightly different early-repayment run).
#   - A couple of structural cleanups (blank cells added/removed, a couple
#     of stray cell styles straightened out to match their neighbours).
#   - Selection / active-sheet bookkeeping: Summary becomes the active tab
#     (it was Transactions before).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsInput   = $wb.Worksheets.Item("Input")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsRepay   = $wb.Worksheets.Item("Repayment Schedule")
$wsTrans   = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------------
# 1. Summary sheet - updated totals
# ---------------------------------------------------------------------------
$wsSummary.Range("B2").Value = 790.76
$wsSummary.Range("E2").Value = 9209.24
$wsSummary.Range("F2").Value = 899.9

$wsSummary.Range("A3").Value = 653.76
$wsSummary.Range("B3").Value = 96.96
$wsSummary.Range("E3").Value = 556.79999999999995
$wsSummary.Range("F3").Value = 84.78

# ---------------------------------------------------------------------------
# 2. Repayment Schedule sheet
# ---------------------------------------------------------------------------

# -- row 2: a few previously-empty cells gain matching (empty) formatted
#    cells, and the Waived/Penalties/Fees columns swap which ones carry an
#    explicit 0.
$wsRepay.Range("A2").Copy()
$wsRepay.Range("B2").PasteSpecial(-4122)
$wsRepay.Range("F2").PasteSpecial(-4122)
$wsRepay.Range("O2").PasteSpecial(-4122)

$wsRepay.Range("H2").ClearContents()
$wsRepay.Range("I2").Value = 0
$wsRepay.Range("J2").ClearContents()
$wsRepay.Range("L2").Value = 0

# -- row 3: Paid Date (D3) is cleared and its odd style corrected, and the
#    always-blank E3 cell's stray style is straightened out to match D3/etc.
$wsRepay.Range("D3").ClearContents()
$wsRepay.Range("D2").Copy()
$wsRepay.Range("D3").PasteSpecial(-4122)
$wsRepay.Range("D2").Copy()
$wsRepay.Range("E3").PasteSpecial(-4122)

$wsRepay.Range("F3").Value = 887.72
$wsRepay.Range("G3").Value = 9112.2800000000007
$wsRepay.Range("H3").Value = 96.96
$wsRepay.Range("K3").Value = 984.68
$wsRepay.Range("P3").Value = 96.96

$wsRepay.Range("F4").Value = 802.94
$wsRepay.Range("G4").Value = 8309.34
$wsRepay.Range("H4").Value = 84.78

$wsRepay.Range("F5").Value = 793.86
$wsRepay.Range("G5").Value = 7515.48
$wsRepay.Range("H5").Value = 93.86

$wsRepay.Range("F6").Value = 813.59
$wsRepay.Range("G6").Value = 6701.89
$wsRepay.Range("H6").Value = 74.13

$wsRepay.Range("F7").Value = 819.42
$wsRepay.Range("G7").Value = 5882.47
$wsRepay.Range("H7").Value = 68.3

$wsRepay.Range("F8").Value = 829.7
$wsRepay.Range("G8").Value = 5052.7700000000004
$wsRepay.Range("H8").Value = 58.02

$wsRepay.Range("F9").Value = 836.22
$wsRepay.Range("G9").Value = 4216.55
$wsRepay.Range("H9").Value = 51.5

$wsRepay.Range("F10").Value = 844.75
$wsRepay.Range("G10").Value = 3371.8
$wsRepay.Range("H10").Value = 42.97

$wsRepay.Range("F11").Value = 854.46
$wsRepay.Range("G11").Value = 2517.34
$wsRepay.Range("H11").Value = 33.26

$wsRepay.Range("F12").Value = 862.06
$wsRepay.Range("G12").Value = 1655.28
$wsRepay.Range("H12").Value = 25.66

$wsRepay.Range("F13").Value = 871.39
$wsRepay.Range("G13").Value = 783.89
$wsRepay.Range("H13").Value = 16.329999999999998

$wsRepay.Range("F14").Value = 783.89
$wsRepay.Range("H14").Value = 7.99
$wsRepay.Range("K14").Value = 791.88
$wsRepay.Range("P14").Value = 791.88

# ---------------------------------------------------------------------------
# 3. Transactions sheet
# ---------------------------------------------------------------------------
$wsTrans.Range("A2").Value = 13
$wsTrans.Range("F2").Value = 790.76
$wsTrans.Range("G2").Value = 96.96
$wsTrans.Range("J2").Value = 9209.24
$wsTrans.Range("A3").Value = 11

# K2:L2 were only ever blank placeholder cells - drop them entirely so the
# sheet's used range shrinks back to A1:J3.
$wsTrans.Range("K2:L2").Clear()

# I1's header cell had picked up a stray "centered" style; straighten it out
# to match the rest of the header row.
$wsTrans.Range("H1").Copy()
$wsTrans.Range("I1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Selection / active-sheet bookkeeping.
#    Order matters: Range.Select() activates its sheet, so the sheet we
#    want active at the end (Summary) must be selected last.
# ---------------------------------------------------------------------------
$wsRepay.Range("D11").Select()
$wsTrans.Range("H2").Select()

$wsSummary.Activate()
$wsSummary.Range("C4").Select()
